$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update DNI values in column A (rows 2-6)
$ws.Range("A2").Value = 11111126
$ws.Range("A3").Value = 11111127
$ws.Range("A4").Value = 11111128
$ws.Range("A5").Value = 11111129
$ws.Range("A6").Value = 11111130

# Update active selection from G8 to C8
$ws.Range("C8").Select()
